# Apply weekly fruit/vegetable price update: rows are re-permuted with
# new source data (row 17 stays, others rotate per the mapping below).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (data pulled from source row's
# ORIGINAL values and written into destination row).
$mapping = @{
    2  = 3
    3  = 4
    4  = 18
    5  = 15
    6  = 8
    7  = 12
    8  = 13
    9  = 14
    10 = 5
    11 = 16
    12 = 6
    13 = 10
    14 = 2
    15 = 7
    16 = 11
    17 = 17
    18 = 9
}

# Columns that actually change per the diff: D, H, I, J, K, L, M, N, P, Q
$cols = @(4, 8, 9, 10, 11, 12, 13, 14, 16, 17)

$orig = @{}
foreach ($r in 2..18) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $orig[$r] = $rowData
}

foreach ($destRow in 2..18) {
    $srcRow = $mapping[$destRow]
    if ($srcRow -eq $destRow) {
        continue
    }
    $srcData = $orig[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value2 = $srcData[$c]
    }
}
